$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.848.36"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "3.717.73"
$ws.Range("E3").Value = "  -2.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("D7").Value = "3.705.24"
$ws.Range("E7").Value = "  -2.63%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000243"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "4.328.07"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").Value = "3.698.50"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "67.680.80"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "489.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000142"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.62%  "
$ws.Range("D34").Value = "3.851.43"
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.109"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("D36").Value = "3.654.82"
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.996"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "429.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0353"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Value = "2.752.40"
$ws.Range("E51").Value = "  -3.96%  "
